$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the photo-attribution related image link / name cells (O2:O8, R2),
# and remove the attribution-name helper column (Q2:Q8) entirely.
$ws.Range("O2:O8").ClearContents()
$ws.Range("R2").ClearContents()
$ws.Range("Q2:Q8").ClearContents()

# Update the selected cell shown when the workbook is opened.
$ws.Range("R3").Select()
